# KAL.xlsx manual upload 2021/07/2 8:00
# - CargoEtsngName ("МЕДЬ") renamed to "КОНЦЕНТР МЕД" for every data row.
# - A handful of CarAmount (column B) values were updated from 0 to the
#   actual shipped volume for specific dates.
# - The active window was scrolled down and a different cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the cargo name used throughout column E (CargoEtsngName).
#    Every populated data row (2..179) currently holds "МЕДЬ"; it becomes
#    "КОНЦЕНТР МЕД" after the edit.
# ---------------------------------------------------------------------------
$lastRow = 179
$ws.Range("E2:E$lastRow").Value = "КОНЦЕНТР МЕД"

# ---------------------------------------------------------------------------
# 2) Update the CarAmount (column B) values that changed from 0 to a
#    non-zero shipment quantity on specific rows.
# ---------------------------------------------------------------------------
$carAmountUpdates = @{
    33 = 30; 34 = 30; 35 = 30; 36 = 30; 37 = 30; 38 = 30; 39 = 30
    41 = 30; 42 = 30; 43 = 30
    45 = 30; 46 = 30
    48 = 30; 49 = 30
    51 = 30; 52 = 30; 53 = 30
    55 = 30; 56 = 30
    58 = 30; 59 = 30; 60 = 30
    61 = 30; 62 = 30; 63 = 30; 64 = 30; 65 = 30
    66 = 34; 67 = 34
    69 = 44; 70 = 44
    86 = 50; 87 = 30; 88 = 50; 89 = 50; 90 = 50
    91 = 16; 93 = 16; 95 = 16; 97 = 16
    102 = 12
    117 = 20
    123 = 21; 125 = 21; 127 = 21; 129 = 21; 133 = 21; 136 = 21; 139 = 21; 143 = 21
    146 = 25
    154 = 20; 157 = 20; 160 = 20; 163 = 20; 166 = 20; 168 = 20
}

foreach ($row in $carAmountUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $carAmountUpdates[$row]
}

# ---------------------------------------------------------------------------
# 3) Scroll the view down and move the selection (best-effort: mirrors the
#    author's new cursor position/viewport from the saved sheetView).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 157
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H164").Select()
